$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while preventing Excel from
# auto-converting number-like strings (e.g. "489.73") into real numbers,
# and without leaving a new number-format style attached to the cell.
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = '57.192.72'
$ws.Range("E2").Value = '  +0.89%  '

$ws.Range("D3").Value = '2.430.02'
$ws.Range("E3").Value = '  -1.48%  '

$ws.Range("E4").Value = '  +0.00%  '

Set-TextValue $ws.Range("D5") '489.73'
$ws.Range("E5").Value = '  +0.12%  '

Set-TextValue $ws.Range("D6") '154.55'
$ws.Range("E6").Value = '  +2.07%  '

$ws.Range("E7").Value = '  +19.94%  '

$ws.Range("E8").Value = '  -0.13%  '

$ws.Range("D9").Value = '2.446.95'
$ws.Range("E9").Value = '  -1.12%  '

Set-TextValue $ws.Range("D10") '6.19'
$ws.Range("E10").Value = '  +8.68%  '

Set-TextValue $ws.Range("D11") '0.100'
$ws.Range("E11").Value = '  +0.95%  '

Set-TextValue $ws.Range("D12") '0.334'
$ws.Range("E12").Value = '  +0.08%  '

$ws.Range("E13").Value = '  +1.18%  '

$ws.Range("D14").Value = '2.848.83'
$ws.Range("E14").Value = '  -1.87%  '

$ws.Range("D15").Value = '57.162.99'
$ws.Range("E15").Value = '  +0.21%  '

Set-TextValue $ws.Range("D16") '20.63'
$ws.Range("E16").Value = '  -1.57%  '

$ws.Range("E17").Value = '  -2.31%  '

$ws.Range("D18").Value = '2.439.96'
$ws.Range("E18").Value = '  -1.36%  '

Set-TextValue $ws.Range("D19") '4.63'
$ws.Range("E19").Value = '  +1.91%  '

Set-TextValue $ws.Range("D20") '324.75'
$ws.Range("E20").Value = '  +1.56%  '

Set-TextValue $ws.Range("D21") '10.02'
$ws.Range("E21").Value = '  -1.09%  '

Set-TextValue $ws.Range("D22") '0.998'
$ws.Range("E22").Value = '  +0.04%  '

Set-TextValue $ws.Range("D23") '5.92'
$ws.Range("E23").Value = '  +1.49%  '

Set-TextValue $ws.Range("D24") '57.88'
$ws.Range("E24").Value = '  -0.17%  '

$ws.Range("B25").Value = 'Binance-PegBSC-USD'
$ws.Range("C25").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range("D25") '0.998'
$ws.Range("E25").Value = '  -0.39%  '

$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range("D26") '0.401'
$ws.Range("E26").Value = '  -0.95%  '

$ws.Range("E27").Value = '  -1.39%  '

$ws.Range("D28").Value = '2.532.21'
$ws.Range("E28").Value = '  -2.15%  '

Set-TextValue $ws.Range("D29") '7.27'
$ws.Range("E29").Value = '  -3.65%  '

$ws.Range("D30").Value = '0.0₃0787'
$ws.Range("E30").Value = '  -1.86%  '

Set-TextValue $ws.Range("D32") '150.27'
$ws.Range("E32").Value = '  -0.35%  '

Set-TextValue $ws.Range("D33") '18.72'
$ws.Range("E33").Value = '  +2.64%  '

$ws.Range("E34").Value = '  +0.77%  '

$ws.Range("E35").Value = '  +2.01%  '

Set-TextValue $ws.Range("D36") '3.79'
$ws.Range("E36").Value = '  +0.93%  '

$ws.Range("E37").Value = '  -0.85%  '

Set-TextValue $ws.Range("D38") '0.818'

Set-TextValue $ws.Range("D39") '0.102'
$ws.Range("E39").Value = '  +7.58%  '

Set-TextValue $ws.Range("D40") '286.14'
$ws.Range("E40").Value = '  +9.39%  '

Set-TextValue $ws.Range("D41") '33.98'
$ws.Range("E41").Value = '  -0.23%  '

$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D42") '3.53'
$ws.Range("E42").Value = '  +0.89%  '

$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D43") '1.38'
$ws.Range("E43").Value = '  -0.19%  '

Set-TextValue $ws.Range("D44") '0.993'
$ws.Range("E44").Value = '  -0.34%  '

Set-TextValue $ws.Range("D45") '0.604'
$ws.Range("E45").Value = '  -0.50%  '

Set-TextValue $ws.Range("D46") '0.0533'
$ws.Range("E46").Value = '  -4.25%  '

Set-TextValue $ws.Range("D47") '10.24'
$ws.Range("E47").Value = '  +0.27%  '

Set-TextValue $ws.Range("D48") '0.0228'
$ws.Range("E48").Value = '  +0.02%  '

Set-TextValue $ws.Range("D49") '4.58'
$ws.Range("E49").Value = '  -3.81%  '

$ws.Range("D50").Value = '1.894.67'
$ws.Range("E50").Value = '  +2.26%  '

Set-TextValue $ws.Range("D51") '17.57'
$ws.Range("E51").Value = '  -0.77%  '
